# Denominadores_D1_17.xlsx template update
# - Remove the two empty placeholder sheets (Hoja2, Hoja3)
# - Rename/re-order the header row on Hoja1 (Spanish -> English column names,
#   several new ILI/ICU/Deaths/Pneu columns inserted)
# - Disable concurrent (multi-threaded) calculation
# - Re-apply the (bold + yyyy/mm/dd) header style to the shifted columns

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Drop the unused Hoja2 / Hoja3 sheets -------------------------------
foreach ($sheetName in @("Hoja2", "Hoja3")) {
    foreach ($sheet in $wb.Worksheets) {
        if ($sheet.Name -eq $sheetName) {
            $sheet.Delete() | Out-Null
        }
    }
}

# --- Turn off concurrent calculation ------------------------------------
$excel.MultiThreadedCalculation.Enabled = $false

# --- Rebuild the header row on Hoja1 ------------------------------------
$ws = $wb.Worksheets.Item("Hoja1")

$headers = [ordered]@{
    "A1"  = "Hospital"
    "B1"  = "StartDateOfWeek"
    "C1"  = "EW"
    "D1"  = "EpiYear"
    "E1"  = "Age_Group"
    "F1"  = "ILINumFem"
    "G1"  = "ILINumMale"
    "H1"  = "ILINumST"
    "I1"  = "ILINumEmerST"
    "J1"  = "ILIDenoFem"
    "K1"  = "ILIDenoMale"
    "L1"  = "ILIDenoST"
    "M1"  = "HospFem"
    "N1"  = "HospMale"
    "O1"  = "HospST"
    "P1"  = "ICUFem"
    "Q1"  = "ICUMale"
    "R1"  = "ICUST"
    "S1"  = "DeathsFem"
    "T1"  = "DeathsMale"
    "U1"  = "DeathsST"
    "V1"  = "PneuFem"
    "W1"  = "PneuMale"
    "X1"  = "PneuST"
    "Y1"  = "CCSARIFem"
    "Z1"  = "CCSARIMale"
    "AA1" = "CCSARIST"
    "AB1" = "VentFem"
    "AC1" = "VentMale"
    "AD1" = "VentST"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Columns E1:AD1 pick up the (bold + yyyy/mm/dd) header style that used to
# live only on B1 (StartDateOfWeek); A1/C1/D1 keep the plain bold style.
$ws.Range("E1:AD1").NumberFormat = "yyyy/mm/dd"
$ws.Range("E1:AD1").Font.Bold = $true
